$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new snapshot column right after column A (old col B..BT shift right to C..BU)
$ws.Columns.Item(2).Insert()

# Give the newly inserted header cell (B1) the same look as the other timestamp headers
# (bold / bordered / centered) by copying formatting from the neighbouring header cell (C1,
# which used to be B1 before the insert), then stamp in the new snapshot timestamp.
$ws.Cells.Item(1, 3).Copy() | Out-Null
$ws.Cells.Item(1, 2).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Cells.Item(1, 2).Value = "2025-12-26 20:18"

# The new snapshot's prices match the previous (most recent) scrape, which now lives in
# column C after the shift, so duplicate C's values into the new column B for every data row.
$lastRow = $ws.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 2).Value = $ws.Cells.Item($r, 3).Value2
}

# The insert leaves column B at Excel's default width; restore it to match its neighbours.
$ws.Columns.Item(2).ColumnWidth = $ws.Columns.Item(3).ColumnWidth

$excel.CutCopyMode = 0
